$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sort the data rows (A2:C19) alphabetically by city (column A) ---
# Using the worksheet Sort object (like the Data > Sort dialog) so that
# Excel records a sortState/sortCondition in the saved worksheet.
$dataRange = $ws.Range("A1:C19")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A19"))
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Clean up: the old placeholder style that used to sit on every data
# cell no longer applies once the header gets the new look below. ---
$ws.Range("A2:C19").Style = "Normal"

# --- Give the header row (A1:C1) a proper "tidied up" look: the
# built-in "40% - Accent5" cell style, bumped to 12pt and bolded. ---
$header = $ws.Range("A1:C1")
$header.Style = "40% - Accent5"
$headerStyle = $wb.Styles("40% - Accent5")
$headerStyle.Font.Size = 12
$headerStyle.Font.Bold = $true
$header.RowHeight = 16

# --- View tidy-up: zoom to 100% and leave the selection where the
# author left off. ---
$ws.Range("G13").Select()
$excel.ActiveWindow.Zoom = 100
